# Applies the "global ligreto parameters" / "missing" marker change:
#  - Adds a new "<missing>" marker string used instead of null values for
#    join rows where either Source1 or Source2 had no matching record.
#  - "detailed" sheet: for FIRST_NAME-group rows (ID 4) that only had a
#    Source1 value, fill the empty Source2 cell (column E) with "<missing>".
#    For rows (ID 5) that only had a Source2 value, fill the empty Source1
#    cell (column D) with "<missing>".
#  - "interlaced" sheet: same idea, but Source1/Source2 values live in
#    adjacent column pairs, one pair per compared column.
#  - A couple of "ID" columns grow their best-fit width by one character
#    once the new text is present.

$wb = $excel.ActiveWorkbook

function Fill-Missing {
    param(
        $ws,
        [int]$FormatSourceRow,
        [int]$FormatSourceCol,
        [int]$Row,
        [int]$Col
    )
    $src = $ws.Cells.Item($FormatSourceRow, $FormatSourceCol)
    $dst = $ws.Cells.Item($Row, $Col)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats - keep the existing "diff" style (s="4")
    $dst.Value = "<missing>"
}

# ---------------------------------------------------------------------
# Sheet "detailed": columns are B=Column Name, C=ID, D=Source1, E=Source2,
# F=Difference, G=Relative.
# ---------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("detailed")

# ID group 4 (rows 24-30): Source1 (D) present, Source2 (E) was empty.
foreach ($r in 24..30) {
    Fill-Missing $wsDetailed $r 3 $r 5
}

# ID group 5 (rows 31-37): Source2 (E or F) present, Source1 (D) was empty.
foreach ($r in 31..37) {
    Fill-Missing $wsDetailed $r 3 $r 4
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Sheet "interlaced": each compared column has a Source1/Source2 pair of
# adjacent columns (C/D FIRST_NAME, E/F NUM_INT, G/H NUM_FLOAT, ...).
# ---------------------------------------------------------------------
$wsInterlaced = $wb.Worksheets.Item("interlaced")

# Row 6 (ID 4): Source1 present, Source2 missing -> fill D,F,H,J,L,N,P.
foreach ($c in 4, 6, 8, 10, 12, 14, 16) {
    Fill-Missing $wsInterlaced 6 2 6 $c
}

# Row 7 (ID 5): Source2 present, Source1 missing -> fill C,E,G,I,K,M,O.
foreach ($c in 3, 5, 7, 9, 11, 13, 15) {
    Fill-Missing $wsInterlaced 7 2 7 $c
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Best-fit width bump on the columns whose content changed: the ID
# columns ("detailed"!C and "interlaced"!B) and the NUM_INT pair on
# "interlaced" (E/F) grow by roughly one character once "<missing>" is
# present among their values.
# ---------------------------------------------------------------------
$wsDetailed.Columns.Item(3).ColumnWidth = 7.1666666666666667

$wsInterlaced.Columns.Item(2).ColumnWidth = 7.1666666666666667
$wsInterlaced.Columns.Item(5).ColumnWidth = 13.8333333333333333
$wsInterlaced.Columns.Item(6).ColumnWidth = 13.8333333333333333
